# New pdf example, and changed column processing
# Move the second data column (Parcel id without dashes) from column D to
# column R, leaving column A (dashed Parcel id) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 4; $r++) {
    $val = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 18).Value2 = $val
    $ws.Cells.Item($r, 4).Value2 = $null
}

$ws.Range("A1").Select()
